# Weekly refresh of "Haba" price rows: each data row (2-43, except the
# untouched row 20) receives the D/J/K/L/M/O/P values that used to belong
# to a different row in the same column block. Columns A,B,C,E,F,G,H,I,N,Q,R
# are identical across every row and are left untouched.
#
# Because this is a permutation (several disjoint cycles, not just simple
# swaps), we must snapshot every source row's values BEFORE writing any of
# them, otherwise an early write would clobber data a later row still needs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# destinationRow -> sourceRow (the row whose old D/J/K/L/M/O/P values move
# into destinationRow)
$rowMap = @{
    2  = 21
    3  = 12
    4  = 28
    5  = 6
    6  = 11
    7  = 26
    8  = 37
    9  = 22
    10 = 43
    11 = 29
    12 = 16
    13 = 2
    14 = 40
    15 = 42
    16 = 14
    17 = 36
    18 = 41
    19 = 31
    21 = 5
    22 = 10
    23 = 39
    24 = 25
    25 = 8
    26 = 38
    27 = 34
    28 = 35
    29 = 4
    30 = 24
    31 = 7
    32 = 19
    33 = 30
    34 = 23
    35 = 17
    36 = 18
    37 = 32
    38 = 33
    39 = 15
    40 = 3
    41 = 13
    42 = 27
    43 = 9
}

$cols = @("D", "J", "K", "L", "M", "O", "P")

# Snapshot of every row's current values, keyed by row number.
$snapshot = @{}
foreach ($r in $rowMap.Values) {
    if (-not $snapshot.ContainsKey($r)) {
        $rowData = @{}
        foreach ($c in $cols) {
            $rowData[$c] = $ws.Range("$c$r").Value2
        }
        $snapshot[$r] = $rowData
    }
}

# Now write each destination row from the snapshot of its source row.
foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $rowData = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value = $rowData[$c]
    }
}

Write-Output "done"
